# Adding the changes we made on may 9th
#
# The sensor-data sheet gets 6 new rows of readings inserted right after
# the header (pushing the existing data down by 6 rows) and 4 new rows of
# readings appended after what is now the last row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 6 fresh rows right after the header row, shifting the
#        existing data (old rows 2-21) down to rows 8-27. ---
$ws.Rows("2:7").Insert()

# Newly inserted rows come out of Insert() carrying the header row's
# formatting (bold/border); strip that so they match the plain data rows.
$ws.Range("A2:H7").ClearFormats()

# --- 2. Fill the 6 newly inserted rows (new timestamps 0..500) ---
$newTop = @(
  @(0,   "falling", -3.195676267147064, 5.127160429954529, -1.443197593092918, 0.0357356183230876,  0.0120645882561802,  0.1313360333442688),
  @(100, "falling", -3.130342268943787, 5.136516356468201, -1.369547128677368, 0.0088575463742017,  0.0383317954838275,  0.0606283769011497),
  @(200, "falling", -3.034864258766174, 5.101877164840698, -1.325036150217056, -0.0204639863222837, 0.0259617734700441,  0.0542142912745475),
  @(300, "falling", -3.194309616088868, 5.024436473846436, -1.315180826187134, -0.0181732401251792, 0.0203112699091434,  -0.0135917514562606),
  @(400, "falling", -3.382834231853486, 5.098868799209595, -1.453447324037552, -0.0335975885391235, -0.0102319931611418, -0.0829249545931816),
  @(500, "falling", -3.195986032485962, 5.139615774154663, -1.564420849084854, -0.0200058370828628, -0.0244346093386411, -0.0164933614432811)
)

$r = 2
foreach ($row in $newTop) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}

# --- 3. Append 4 more rows of new readings (timestamps 2600..2900) after
#        the last existing data row (now row 27). ---
$newBottom = @(
  @(2600, "falling", 2.14622653722763,  5.513726615905762, -1.284043130278588, 0.0577267669141292,  0.2557998299598694,  0.0554360225796699),
  @(2700, "falling", 2.037818813323974, 5.183717918395995, -1.269947481155395, -0.0160352122038602, -0.030695978552103,  -0.0510072484612464),
  @(2800, "falling", 2.045576536655426, 5.118093979358673, -1.376126399636268, -0.0320704244077205, -0.107512280344963,  -0.04505131021142),
  @(2900, "falling", 2.174056196212769, 5.255697178840638, -1.429987555742264, 0.0154243474826216,  0.1117883399128913,  0.0210748501121997)
)

$r = 28
foreach ($row in $newBottom) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c++
  }
  $r++
}

# --- 4. The "timestamp" (A) and "label" (B) columns are not part of the
#        shifted sensor reading data - they simply follow the row's
#        position (100ms increments / constant "falling" label). Re-stamp
#        them for every data row so the values inherited from the
#        row-insert shift (which moved the whole row, A/B included) line
#        back up with the row they now sit in. ---
for ($r = 2; $r -le 31; $r++) {
  $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
  $ws.Cells.Item($r, 2).Value = "falling"
}
